$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "RM 232" (row 26 in the original sheet)
$ws.Rows(26).Delete()

# After the above deletion, the "SC 92" row shifts from row 28 to row 27; delete it too
$ws.Rows(27).Delete()

# Now apply remaining cell-level value changes using final (post-deletion) row numbers
$ws.Range("E6").Value = -5.7
$ws.Range("E8").ClearContents()
$ws.Range("E12").Value = -5.3
$ws.Range("E14").ClearContents()
$ws.Range("E17").Value = -7.3
$ws.Range("E18").Value = -8.5
$ws.Range("E19").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("E23").Value = -7

$ws.Range("D27").Value = -14.6
$ws.Range("E27").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("D32").ClearContents()
